$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("A2").Value = "052/FKIH BEN SALEH"
$ws.Range("B2").Value = "Point de vente"
$ws.Range("C2").Value = "IB19558"
$ws.Range("D2").Value = "ZERNAKH ABDELLAH"
$ws.Range("E2").Value = "oui"
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = $true
$ws.Range("J2").Value = 0
$ws.Range("O2").Value = 11000

# --- Row 3 ---
$ws.Range("A3").Value = "052/FKIH BEN SALEH"
$ws.Range("B3").Value = "Point de vente"
$ws.Range("C3").Value = "IB43905"
$ws.Range("D3").Value = "NHILA BELGACEM"
$ws.Range("E3").Value = "oui"
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = $false
$ws.Range("J3").Value = 0
$ws.Range("O3").Value = 0

# --- Row 4 ---
$ws.Range("A4").Value = "905/TADLA OUARDIGHA ZAYANE"
$ws.Range("C4").Value = "Q251990"
$ws.Range("D4").Value = "NOUBAIL MOUNTASSIR"
$ws.Range("G4").Value = 10
$ws.Range("H4").Value = $true
$ws.Range("J4").Value = 675
$ws.Range("O4").Value = 6075

# --- Row 5 ---
$ws.Range("A5").Value = "905/TADLA OUARDIGHA ZAYANE"
$ws.Range("B5").Value = "Direction régionale"
$ws.Range("C5").Value = "IR801997"
$ws.Range("D5").Value = "NOUBAIL MOHAMMED"
$ws.Range("H5").Value = $true
$ws.Range("J5").Value = 675
$ws.Range("O5").Value = 6075

# --- Row 6 (totals) ---
$ws.Range("H6").Value = 3
$ws.Range("J6").Value = 1350
$ws.Range("O6").Value = 23150
